$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.6541118602811586
$ws.Range("E2").Value = 0.6541118602811586

$ws.Range("D3").Value = 0.03323614730896127
$ws.Range("E3").Value = 0.03323614730896127

$ws.Range("D4").Value = 0.9730133167169314
$ws.Range("E4").Value = 0.9730133167169314

$ws.Range("D5").Value = 0.09884951202197124
$ws.Range("E5").Value = 0.09884951202197124

$ws.Range("D6").Value = 0.3067804907064885
$ws.Range("E6").Value = 0.3067804907064885

$ws.Range("D7").Value = 0.9999892282737609
$ws.Range("E7").Value = 0.0000107717262390894

$ws.Range("D8").Value = 0.9735798597716366
$ws.Range("E8").Value = 0.02642014022836336

$ws.Range("D9").Value = 0.9778197031698652
$ws.Range("E9").Value = 0.02218029683013478

$ws.Range("D10").Value = 0.998930254650726
$ws.Range("E10").Value = 0.001069745349273954

$ws.Range("D11").Value = 0.9763926654641428
$ws.Range("E11").Value = 0.02360733453585717
$ws.Range("F11").Value = 0.5252521634101868

$ws.Range("D12").Value = 0.6817748442996882
$ws.Range("E12").Value = 0.6817748442996882

$ws.Range("D13").Value = 0.00109120342154986
$ws.Range("E13").Value = 0.00109120342154986

$ws.Range("D14").Value = 0.9646799378456709
$ws.Range("E14").Value = 0.9646799378456709

$ws.Range("D15").Value = 0.0001425754764067132
$ws.Range("E15").Value = 0.0001425754764067132

$ws.Range("D16").Value = 0.1364974265987901
$ws.Range("E16").Value = 0.1364974265987901

$ws.Range("D17").Value = 0.9999855817346953
$ws.Range("E17").Value = 0.00001441826530468315

$ws.Range("D18").Value = 0.9910876821438145
$ws.Range("E18").Value = 0.008912317856185514

$ws.Range("D19").Value = 0.9868300976855224
$ws.Range("E19").Value = 0.01316990231447757

$ws.Range("D20").Value = 0.9993609739937203
$ws.Range("E20").Value = 0.0006390260062797415

$ws.Range("D21").Value = 0.9991173853432029
$ws.Range("E21").Value = 0.000882614656797065
$ws.Range("F21").Value = 0.466003954410553
